$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on Price (D) and Volume (E) columns for the data rows
# so that numeric-looking strings (e.g. "1.001", "7.337") are preserved
# exactly as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.017.78'
$ws.Range("E2").Value = '  -3.47%  '
$ws.Range("D3").Value = '1.600.50'
$ws.Range("E3").Value = '  -2.34%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '301.14'
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("D7").Value = '0.3781'
$ws.Range("E7").Value = '  -2.29%  '
$ws.Range("D8").Value = '0.3638'
$ws.Range("E8").Value = '  -4.38%  '
$ws.Range("D9").Value = '49.89'
$ws.Range("E9").Value = '  -1.14%  '
$ws.Range("E10").Value = '  -4.81%  '
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '0.08126'
$ws.Range("E12").Value = '  -2.92%  '
$ws.Range("D13").Value = '22.58'
$ws.Range("E13").Value = '  -4.32%  '
$ws.Range("E14").Value = '  -5.11%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.337'
$ws.Range("E15").Value = '  -6.26%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").Value = '0.00001243'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("D17").Value = '1.605.58'
$ws.Range("E17").Value = '  -2.07%  '
$ws.Range("D18").Value = '91.78'
$ws.Range("E18").Value = '  -1.66%  '
$ws.Range("D19").Value = '0.06814'
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").Value = '18.21'
$ws.Range("E20").Value = '  -5.96%  '
$ws.Range("D21").Value = '6.550'
$ws.Range("E21").Value = '  -4.34%  '
$ws.Range("B22").Value = 'BitDAO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range("D22").Value = '0.5572'
$ws.Range("E22").Value = '  -5.58%  '
$ws.Range("B23").Value = 'Dai'
$ws.Range("C23").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '13.09'
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D25").Value = '23.026.31'
$ws.Range("E25").Value = '  -3.46%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").Value = '2.355'
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").Value = '2.821'
$ws.Range("E27").Value = '  -1.57%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '21.03'
$ws.Range("E28").Value = '  -3.34%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").Value = '150.40'
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").Value = '5.242'
$ws.Range("E30").Value = '  -5.13%  '
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").Value = '133.62'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D32").Value = '2.290'
$ws.Range("E32").Value = '  -8.20%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '6.833'
$ws.Range("E33").Value = '  -12.04%  '
$ws.Range("B34").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C34").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D34").Value = '1.785.58'
$ws.Range("E34").Value = '  -1.95%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.9655'
$ws.Range("E35").Value = '  -1.17%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.07581'
$ws.Range("E36").Value = '  -4.56%  '
$ws.Range("B37").Value = 'FraxShare'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D37").Value = '10.33'
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = '6.266'
$ws.Range("E38").Value = '  -4.92%  '
$ws.Range("D39").Value = '0.02703'
$ws.Range("E39").Value = '  -6.51%  '
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.2532'
$ws.Range("E40").Value = '  -4.46%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.08887'
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '1.370'
$ws.Range("E42").Value = '  -2.79%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.7026'
$ws.Range("E43").Value = '  -5.90%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '12.44'
$ws.Range("E44").Value = '  -6.05%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '15.17'
$ws.Range("E45").Value = '  -8.00%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = '0.6636'
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.308'
$ws.Range("E48").Value = '  -4.09%  '
$ws.Range("B49").Value = 'PancakeSwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D49").Value = '3.989'
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '132.30'
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("D51").Value = '0.07922'
